{"js": "// Delete the \"\u5831\u544a\u65e5: 2024 \u5e74 1 \u6708 22 \u65e5\" (report date) paragraph and the\n// blank (single-space) paragraph that immediately follows it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"\u5831\u544a\u65e5: 2024 \u5e74 1 \u6708 22 \u65e5\";\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  // The paragraph right after it is the lone-space paragraph that should\n  // be removed together with the date paragraph.\n  const dateParagraph = paragraphs.items[targetIndex];\n  const spaceParagraph =\n    targetIndex + 1 < paragraphs.items.length\n      ? paragraphs.items[targetIndex + 1]\n      : null;\n\n  if (spaceParagraph && spaceParagraph.text.trim() === \"\") {\n    spaceParagraph.delete();\n  }\n  dateParagraph.delete();\n\n  await context.sync();\n}\n", "ps1": "# Delete the \"\u5831\u544a\u65e5: 2024 \u5e74 1 \u6708 22 \u65e5\" (report date) paragraph and the\n# blank (single-space) paragraph that immediately follows it.\n$d = $word.ActiveDocument\n\n$targetText = \"\u5831\u544a\u65e5: 2024 \u5e74 1 \u6708 22 \u65e5\"\n$foundIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq $targetText) {\n        $foundIndex = $i\n        break\n    }\n}\n\nif ($foundIndex -ge 1) {\n    if ($foundIndex + 1 -le $d.Paragraphs.Count) {\n        $nextPara = $d.Paragraphs.Item($foundIndex + 1)\n        if ($nextPara.Range.Text.Trim() -eq \"\") {\n            $nextPara.Range.Delete()\n        }\n    }\n    $d.Paragraphs.Item($foundIndex).Range.Delete()\n}\n"}
